# Refresh the cryptocurrency price list (column D) and 1h change
# percentages (column E) for rows 2-51 with the latest scraped values.
#
# Column D holds price strings that are stored as TEXT in the sheet
# (e.g. "58.083.73" uses "." as a thousands separator, which is not a
# valid Excel number). A leading apostrophe is Excel's standard text
# qualifier: it forces Range.Value assignment to keep the string as text
# (instead of auto-converting values like "528.07" to a number) while the
# apostrophe itself is not stored in the cell's value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.121.95"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "'3.120.17"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'528.07"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'142.46"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'3.117.91"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "'3.654.22"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("D15").Value = "'25.60"
$ws.Range("E15").Value = "  -4.30%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "'58.153.52"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "'3.118.67"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "'6.12"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "'12.79"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'7.98"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").Value = "'342.48"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'67.60"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'0.0₃0925"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'7.33"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").Value = "'1.88"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "'21.05"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("D35").Value = "'158.36"
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("D37").Value = "'6.21"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").Value = "'26.42"
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("E40").Value = "  +12.85%  "
$ws.Range("D41").Value = "'0.0667"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").Value = "'0.692"
$ws.Range("E43").Value = "  +4.28%  "
$ws.Range("D44").Value = "'3.157.47"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").Value = "'36.63"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").Value = "'2.274.45"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +4.50%  "
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").Value = "'20.68"
$ws.Range("E51").Value = "  -0.09%  "
